$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '26.848.16'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +0.42%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.641.02'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("E5").Value = '  +0.62%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '0.497'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.61%  '

$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("E8").Value = '  -0.45%  '

$ws.Range("E10").Value = '  +0.60%  '

$ws.Range("E11").Value = '  +0.37%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '1.870.41'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.02%  '

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.642.47'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.01%  '

$ws.Range("E14").Value = '  -0.30%  '

$ws.Range("E15").Value = '  +0.10%  '

$ws.Range("E16").Value = '  +1.24%  '

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '26.854.02'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +0.44%  '

$ws.Range("E18").Value = '  -0.71%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '214.71'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +0.67%  '

$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("E21").Value = '  -0.01%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '6.59'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +5.53%  '

$ws.Range("E23").Value = '  -3.09%  '

$ws.Range("E24").Value = '  -1.40%  '

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '147.50'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.33%  '

$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("E27").Value = '  -0.20%  '

$ws.Range("E28").Value = '  +1.42%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '15.73'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +0.60%  '

$ws.Range("E30").Value = '  +0.03%  '

$ws.Range("E31").Value = '  +1.05%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.36'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '

$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.55'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +1.35%  '

$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.277.55'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -1.11%  '

$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("E37").Value = '  -1.51%  '

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.530'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -0.55%  '

$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("E40").Value = '  -0.14%  '

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.804'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("E42").Value = '  +0.04%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.780.88'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.60%  '

$ws.Range("E44").Value = '  -6.01%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '92.46'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +1.39%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '61.04'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.35%  '

$ws.Range("E47").Value = '  -0.55%  '

$ws.Range("E48").Value = '  -1.74%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.0967'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.62%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '7.56'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -0.94%  '

$ws.Range("E51").Value = '  -0.06%  '
